# Apply cell updates from the "Updated symbol list" commit (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''278.04'
$ws.Range("E2").Value = '''0.49%'
$ws.Range("G2").Value = '''11'
$ws.Range("D3").Value = '''27.21'
$ws.Range("E3").Value = '''1.78%'
$ws.Range("G3").Value = '''11'
$ws.Range("D4").Value = '''4.868'
$ws.Range("G4").Value = '''11'
$ws.Range("D5").Value = '''0.06431'
$ws.Range("E5").Value = '''1.63%'
$ws.Range("G5").Value = '''11'
$ws.Range("D6").Value = '''7.027'
$ws.Range("E6").Value = '''1.41%'
$ws.Range("G6").Value = '''11'
$ws.Range("D7").Value = '''1.194'
$ws.Range("E7").Value = '''-6.54%'
$ws.Range("G7").Value = '''11'
$ws.Range("D8").Value = '''0.8850'
$ws.Range("E8").Value = '''0.96%'
$ws.Range("G8").Value = '''11'
$ws.Range("D9").Value = '''0.1548'
$ws.Range("E9").Value = '''-0.51%'
$ws.Range("G9").Value = '''11'
$ws.Range("D10").Value = '''0.05122'
$ws.Range("E10").Value = '''1.59%'
$ws.Range("G10").Value = '''11'
$ws.Range("D11").Value = '''0.07508'
$ws.Range("E11").Value = '''0.48%'
$ws.Range("G11").Value = '''11'
$ws.Range("D12").Value = '''0.02891'
$ws.Range("E12").Value = '''-1.80%'
$ws.Range("G12").Value = '''11'
$ws.Range("D13").Value = '''0.08970'
$ws.Range("E13").Value = '''-0.97%'
$ws.Range("G13").Value = '''11'
$ws.Range("D14").Value = '''0.001562'
$ws.Range("E14").Value = '''-0.70%'
$ws.Range("G14").Value = '''11'
$ws.Range("D15").Value = '''0.0006400'
$ws.Range("E15").Value = '''0.96%'
$ws.Range("G15").Value = '''11'
$ws.Range("D16").Value = '''0.006092'
$ws.Range("E16").Value = '''1.07%'
$ws.Range("G16").Value = '''11'
$ws.Range("D17").Value = '''3.477'
$ws.Range("E17").Value = '''0.80%'
$ws.Range("G17").Value = '''11'
$ws.Range("E18").Value = '''-0.44%'
$ws.Range("G18").Value = '''11'
$ws.Range("G19").Value = '''11'
$ws.Range("G20").Value = '''11'
$ws.Range("D21").Value = '''0.1341'
$ws.Range("E21").Value = '''0.83%'
$ws.Range("G21").Value = '''11'
$ws.Range("D22").Value = '''3.908'
$ws.Range("E22").Value = '''0.14%'
$ws.Range("G22").Value = '''11'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.1519'
$ws.Range("E23").Value = '''10.03%'
$ws.Range("G23").Value = '''11'
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").Value = '''0.04411'
$ws.Range("E24").Value = '''0.99%'
$ws.Range("G24").Value = '''11'
$ws.Range("D25").Value = '''0.001177'
$ws.Range("E25").Value = '''0.30%'
$ws.Range("G25").Value = '''11'
$ws.Range("D26").Value = '''0.003877'
$ws.Range("E26").Value = '''-7.93%'
$ws.Range("G26").Value = '''11'
$ws.Range("G27").Value = '''11'
$ws.Range("D28").Value = '''0.0001181'
$ws.Range("G28").Value = '''11'
$ws.Range("E29").Value = '''1.78%'
$ws.Range("G29").Value = '''11'
$ws.Range("G30").Value = '''11'
$ws.Range("G31").Value = '''11'
$ws.Range("G32").Value = '''11'
$ws.Range("G33").Value = '''11'
$ws.Range("G34").Value = '''11'
$ws.Range("G35").Value = '''11'
$ws.Range("G36").Value = '''11'
$ws.Range("G37").Value = '''11'
$ws.Range("G38").Value = '''11'
$ws.Range("G39").Value = '''11'
$ws.Range("D40").Value = '''0.04123'
$ws.Range("E40").Value = '''0.33%'
$ws.Range("G40").Value = '''11'
$ws.Range("D41").Value = '''0.006806'
$ws.Range("E41").Value = '''-2.87%'
$ws.Range("G41").Value = '''11'
$ws.Range("E42").Value = '''0.04%'
$ws.Range("G42").Value = '''11'
$ws.Range("E43").Value = '''-13.48%'
$ws.Range("G43").Value = '''11'
$ws.Range("D44").Value = '''0.01164'
$ws.Range("E44").Value = '''1.34%'
$ws.Range("G44").Value = '''11'
$ws.Range("D45").Value = '''0.00005321'
$ws.Range("E45").Value = '''0.37%'
$ws.Range("G45").Value = '''11'
$ws.Range("D46").Value = '''1.682'
$ws.Range("E46").Value = '''12.94%'
$ws.Range("G46").Value = '''11'
$ws.Range("E47").Value = '''-7.38%'
$ws.Range("G47").Value = '''11'
$ws.Range("G48").Value = '''11'
$ws.Range("G49").Value = '''11'
$ws.Range("G50").Value = '''11'
$ws.Range("G51").Value = '''11'

Write-Host "Applied 113 cell updates."
